$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "severity" reference rows (46-50): prefix labels with
# their numeric rank, e.g. "Tidak Penting " -> "0- Tidak Penting" ---
$ws.Range("C46").Value = "0- Not Important"
$ws.Range("D46").Value = "0- Tidak Penting"

$ws.Range("C47").Value = "1- Critical"
$ws.Range("D47").Value = "1- Kritikal"

$ws.Range("C48").Value = "2- Important"
$ws.Range("D48").Value = "2- Penting"

$ws.Range("C49").Value = "3- Medium"
$ws.Range("D49").Value = "3- Sederhana"

$ws.Range("C50").Value = "4- Low"
$ws.Range("D50").Value = "4- Rendah"

# --- Append a new "loaner_type" reference table (rows 63-65) ---
$ws.Range("A63").Value = "loaner_type"
$ws.Range("B63").Value = 1
$ws.Range("C63").Value = "Day"
$ws.Range("D63").Value = "Hari"

$ws.Range("A64").Value = "loaner_type"
$ws.Range("B64").Value = 2
$ws.Range("C64").Value = "Week"
$ws.Range("D64").Value = "Minggu"

$ws.Range("A65").Value = "loaner_type"
$ws.Range("B65").Value = 3
$ws.Range("C65").Value = "Month"
$ws.Range("D65").Value = "Bulan"

# --- Update the window scroll position / selection to match where the
# author was working when the file was saved ---
$ws.Range("A49").Select()
$excel.ActiveWindow.ScrollRow = 49
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C65").Select()
